# Fix pre/post cut column: add TableID markers for the "cflux" block and a
# new "reflectance" table (time/ndvi/notes/pre_post_cut) to the data
# dictionary sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("table_description")

# Fill in the variable names for the new "reflectance" table rows first.
$ws.Range("B57").Value = "ndvi"
$ws.Range("B58").Value = "notes"
$ws.Range("B59").Value = "pre_post_cut"

# Then the first new row's description / units.
$ws.Range("C56").Value = "Time of sampling"
$ws.Range("D56").Value = "hh:mm"

# Then the remaining descriptions.
$ws.Range("C57").Value = "NDVI value"
$ws.Range("C58").Value = "Notes"
$ws.Range("C59").Value = "Measurment was taken before or after the cut"

# Mark the start of the existing "cflux" table (row 32 was missing its
# TableID label in column A), then the new "reflectance" table's TableID.
$ws.Range("A32").Value = "cflux"
$ws.Range("A56").Value = "reflectance"

# Fill in the remaining cells (these all reuse existing shared strings).
$ws.Range("B56").Value = "time"
$ws.Range("E56").Value = "defined"
$ws.Range("D57").Value = "percentage"
$ws.Range("E57").Value = "measured"
$ws.Range("D59").Value = "pre or post"
$ws.Range("E59").Value = "recorded"

# Restore the view to match the saved workbook state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 12
$ws.Range("D26").Select()
